$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 117999
$ws.Range("J63").Value = 117999
$ws.Range("L63").Value = 117999
$ws.Range("N63").Value = -119247
$ws.Range("H66").Value = 117999
$ws.Range("J66").Value = 117999
$ws.Range("L66").Value = 353997
$ws.Range("N66").Value = -360237
$ws.Range("H125").Value = 4360.9165
$ws.Range("I125").Value = 3099.6667
$ws.Range("J125").Value = 4781.3335
$ws.Range("K125").Value = 27897.0003
$ws.Range("L125").Value = 43032.0015
$ws.Range("M125").Value = -25437.0003
$ws.Range("N125").Value = -47952.0015
$ws.Range("H137").Value = 1177092.2
$ws.Range("I137").Value = 70779.664
$ws.Range("J137").Value = 1545863.1
$ws.Range("K137").Value = 212338.992
$ws.Range("L137").Value = 4637589.300000001
$ws.Range("M137").Value = -209788.992
$ws.Range("N137").Value = -4642689.300000001
$ws.Range("H138").Value = 4202.91
$ws.Range("I138").Value = 2026.4546
$ws.Range("J138").Value = 4471.91
$ws.Range("K138").Value = 6079.3638
$ws.Range("L138").Value = 13415.73
$ws.Range("M138").Value = -939.3638000000001
$ws.Range("N138").Value = -23695.73
$ws.Range("H141").Value = 4980.4644
$ws.Range("I141").Value = 4757.5186
$ws.Range("K141").Value = 14272.5558
$ws.Range("M141").Value = -9092.555800000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14022277
$ws.Range("I32").Value = 14198321
$ws.Range("K32").Value = 14198321
$ws.Range("M32").Value = -14198034
$ws.Range("H37").Value = 51295.3
$ws.Range("J37").Value = 68647.5
$ws.Range("L37").Value = 68647.5
$ws.Range("N37").Value = -69193.5
$ws.Range("H74").Value = 2044.9697
$ws.Range("I74").Value = 2044.9697
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2044.9697
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1170.9697
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 2044.9697
$ws.Range("I77").Value = 2044.9697
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 10224.8485
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -5856.8485
$ws.Range("N77").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2036.2593
$ws.Range("I99").Value = 1850.7826
$ws.Range("J99").Value = 3102.75
$ws.Range("K99").Value = 1850.7826
$ws.Range("L99").Value = 3102.75
$ws.Range("M99").Value = -352.7826
$ws.Range("N99").Value = -6098.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4751.9546
$ws.Range("I31").Value = 3997.5
$ws.Range("J31").Value = 4787.881
$ws.Range("K31").Value = 3997.5
$ws.Range("L31").Value = 4787.881
$ws.Range("M31").Value = -3702.5
$ws.Range("N31").Value = -5377.881
$ws.Range("H33").Value = 6174
$ws.Range("I33").Value = 6174
$ws.Range("K33").Value = 6174
$ws.Range("M33").Value = -5795
$ws.Range("H34").Value = 4751.9546
$ws.Range("I34").Value = 3997.5
$ws.Range("J34").Value = 4787.881
$ws.Range("K34").Value = 3997.5
$ws.Range("L34").Value = 4787.881
$ws.Range("M34").Value = -3795.5
$ws.Range("N34").Value = -5191.881
$ws.Range("H122").Value = 3872.5
$ws.Range("I122").Value = 3326.3845
$ws.Range("J122").Value = 4661.3335
$ws.Range("K122").Value = 9979.1535
$ws.Range("L122").Value = 13984.0005
$ws.Range("M122").Value = -7529.1535
$ws.Range("N122").Value = -18884.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 283.875
$ws.Range("I2").Value = 136.5
$ws.Range("J2").Value = 333
$ws.Range("K2").Value = 819
$ws.Range("L2").Value = 1998
$ws.Range("M2").Value = -706
$ws.Range("N2").Value = -2224
$ws.Range("H40").Value = 271.4
$ws.Range("I40").Value = 70
$ws.Range("J40").Value = 573.5
$ws.Range("K40").Value = 280
$ws.Range("L40").Value = 2294
$ws.Range("M40").Value = -211
$ws.Range("N40").Value = -2432
$ws.Range("H68").Value = 3110.7896
$ws.Range("I68").Value = 2785.4285
$ws.Range("J68").Value = 3300.5833
$ws.Range("K68").Value = 8356.2855
$ws.Range("L68").Value = 9901.749899999999
$ws.Range("M68").Value = -7545.2855
$ws.Range("N68").Value = -11523.7499
$ws.Range("H71").Value = 3110.7896
$ws.Range("I71").Value = 2785.4285
$ws.Range("J71").Value = 3300.5833
$ws.Range("K71").Value = 25068.8565
$ws.Range("L71").Value = 29705.2497
$ws.Range("M71").Value = -21012.8565
$ws.Range("N71").Value = -37817.2497
$ws.Range("H107").Value = 668.46155
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H122").Value = 996226.7
$ws.Range("I122").Value = 648.5
$ws.Range("J122").Value = 6969696
$ws.Range("K122").Value = 5836.5
$ws.Range("L122").Value = 62727264
$ws.Range("M122").Value = -3386.5
$ws.Range("N122").Value = -62732164
$ws.Range("H132").Value = 1114500.5
$ws.Range("I132").Value = 3000
$ws.Range("J132").Value = 2003700.8
$ws.Range("K132").Value = 27000
$ws.Range("L132").Value = 18033307.2
$ws.Range("M132").Value = -24470
$ws.Range("N132").Value = -18038367.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3443.5715
$ws.Range("I22").Value = 2600.8
$ws.Range("J22").Value = 3911.7778
$ws.Range("K22").Value = 2600.8
$ws.Range("L22").Value = 3911.7778
$ws.Range("M22").Value = -2305.8
$ws.Range("N22").Value = -4501.7778
$ws.Range("H27").Value = 3443.5715
$ws.Range("I27").Value = 2600.8
$ws.Range("J27").Value = 3911.7778
$ws.Range("K27").Value = 2600.8
$ws.Range("L27").Value = 3911.7778
$ws.Range("M27").Value = -2493.8
$ws.Range("N27").Value = -4125.7778
$ws.Range("H132").Value = 5423.4546
$ws.Range("I132").Value = 5423.4546
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 16270.3638
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -13740.3638
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 6467.5557
$ws.Range("I136").Value = 5569.9165
$ws.Range("K136").Value = 16709.7495
$ws.Range("M136").Value = -14159.7495

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5380.7144
$ws.Range("I81").Value = 6666
$ws.Range("J81").Value = 5166.5
$ws.Range("K81").Value = 13332
$ws.Range("L81").Value = 10333
$ws.Range("M81").Value = -12271
$ws.Range("N81").Value = -12455
$ws.Range("H84").Value = 5380.7144
$ws.Range("I84").Value = 6666
$ws.Range("J84").Value = 5166.5
$ws.Range("K84").Value = 66660
$ws.Range("L84").Value = 51665
$ws.Range("M84").Value = -61356
$ws.Range("N84").Value = -62273
$ws.Range("H122").Value = 4142.3687
$ws.Range("I122").Value = 4110.161
$ws.Range("J122").Value = 4285
$ws.Range("K122").Value = 12330.483
$ws.Range("L122").Value = 12855
$ws.Range("M122").Value = -9880.483
$ws.Range("N122").Value = -17755
$ws.Range("H132").Value = 2460.1785
$ws.Range("I132").Value = 2173.1
$ws.Range("J132").Value = 3177.875
$ws.Range("K132").Value = 6519.299999999999
$ws.Range("L132").Value = 9533.625
$ws.Range("M132").Value = -3989.299999999999
$ws.Range("N132").Value = -14593.625
